$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert two new rows at the top of the question list (rows 3-4), pushing
# the existing "main_stem" row (and everything after it) down by two rows.
$ws.Rows("3:4").Insert()

# New row 3: integer field "plot" / "Plot"
$ws.Range("C3").Value = "integer"
$ws.Range("E3").Value = "plot"
$ws.Range("F3").Value = "Plot"

# New row 4: integer field "tag" / "Tag"
$ws.Range("C4").Value = "integer"
$ws.Range("E4").Value = "tag"
$ws.Range("F4").Value = "Tag"

# Fix the mislabeled "support_percentage" prompt text (now row 11 after the
# two-row insert) so its display prompt reads "Support percentage" instead
# of repeating the field's machine name.
$ws.Range("F11").Value = "Support percentage"

# The "survey" sheet becomes the active tab/selection in the saved workbook.
$ws.Activate()
$ws.Range("C4:F4").Select()
